$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 12, shifting the existing rows 12:17 down to 13:18
$ws.Rows("12:12").Insert()

# Populate the newly inserted row 12 with the new weekly entry
$ws.Cells.Item(12, 1).Value = 1
$ws.Cells.Item(12, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(12, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(12, 4).Value = 44596
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(13, 4).NumberFormat
$ws.Cells.Item(12, 5).Value = 15
$ws.Cells.Item(12, 6).Value = "Fruta"
$ws.Cells.Item(12, 7).Value = 100103
$ws.Cells.Item(12, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(12, 9).Value = 100103002
$ws.Cells.Item(12, 10).Value = "Ciruela"
$ws.Cells.Item(12, 11).Value = "Black Amber"
$ws.Cells.Item(12, 12).Value = "Segunda"
$ws.Cells.Item(12, 13).Value = 250
$ws.Cells.Item(12, 14).Value = 15000
$ws.Cells.Item(12, 15).Value = 16000
$ws.Cells.Item(12, 16).Value = 15500
$ws.Cells.Item(12, 17).Value = "`$/caja 18 kilos granel"
$ws.Cells.Item(12, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(12, 19).Value = 861
$ws.Cells.Item(12, 20).Value = 18
